$d = $word.ActiveDocument

function Get-ParagraphContaining($text) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $ok = $f.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $text"
    }
    return $f.Parent.Paragraphs(1)
}

# ---------------------------------------------------------------------------
# 1) First paragraph: "Munashe Mugonda, Cody Burker"
#    Remove the spell-check proofErr wrappers and merge the four runs that
#    spell out the names into a single run, leaving the visible text as-is.
# ---------------------------------------------------------------------------
$pName = Get-ParagraphContaining("Munashe")
$rName = $pName.Range
$xmlName = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="533172D5" w14:textId="77777777" w:rsidR="00525C90" w:rsidRDefault="00525C90"><w:r><w:t>Munashe Mugonda, Cody Burker</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rName.InsertXML($xmlName) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Materials" paragraph: drop the stray _GoBack bookmark pair that follows
#    the word "Materials".
# ---------------------------------------------------------------------------
$pMaterials = Get-ParagraphContaining("Materials")
$rMaterials = $pMaterials.Range
$xmlMaterials = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="5E177BC0" w14:textId="30EAB55B" w:rsidR="00140E68" w:rsidRDefault="00140E68" w:rsidP="003C5067"><w:r><w:t>Materials</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rMaterials.InsertXML($xmlMaterials) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Research different methods to control Arduino from desktop": merge the
#    two split runs into a single run (plain text Find/Replace merges
#    same-formatted runs cleanly here, since there is no other markup in the
#    way).
# ---------------------------------------------------------------------------
$searchText = "Research different met" + "hods to control Arduino from desktop"
$d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) The empty bullet right after it gains the text "Install Git" + "ub",
#    with the _GoBack bookmark sitting between the two runs (reproducing the
#    exact run/bookmark split recorded in the source revision).
# ---------------------------------------------------------------------------
$pSteps = Get-ParagraphContaining("Research different methods to control Arduino from desktop")
$pGit = $pSteps.Next()
$rGit = $pGit.Range
$xmlGit = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="0549A86E" w14:textId="77777777" w:rsidR="009A771E" w:rsidRDefault="009A771E" w:rsidP="007054C2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Install Git</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ub</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rGit.InsertXML($xmlGit) | Out-Null
